{"js": "// Apply the text replacements described by the diff:\n// - update the date line\n// - update each two-digit multiplication problem in the table\n\nconst replacements = [\n  [\"2026-02-23 Monday\", \"2026-02-24 Tuesday\"],\n  [\"99\u00d798=\", \"26\u00d738=\"],\n  [\"43\u00d722=\", \"97\u00d793=\"],\n  [\"67\u00d789=\", \"66\u00d796=\"],\n  [\"51\u00d728=\", \"98\u00d747=\"],\n  [\"18\u00d714=\", \"43\u00d726=\"],\n  [\"47\u00d741=\", \"67\u00d768=\"],\n  [\"60\u00d766=\", \"73\u00d734=\"],\n  [\"69\u00d789=\", \"53\u00d740=\"],\n  [\"17\u00d739=\", \"71\u00d721=\"],\n  [\"91\u00d765=\", \"43\u00d797=\"],\n  [\"41\u00d718=\", \"90\u00d773=\"],\n  [\"35\u00d764=\", \"11\u00d734=\"],\n  [\"19\u00d787=\", \"70\u00d782=\"],\n  [\"58\u00d739=\", \"26\u00d714=\"],\n  [\"21\u00d765=\", \"85\u00d798=\"],\n  [\"43\u00d729=\", \"71\u00d711=\"],\n  [\"70\u00d779=\", \"42\u00d777=\"],\n  [\"73\u00d761=\", \"56\u00d770=\"],\n  [\"40\u00d782=\", \"47\u00d743=\"],\n  [\"92\u00d711=\", \"55\u00d743=\"],\n  [\"96\u00d786=\", \"94\u00d773=\"],\n  [\"70\u00d759=\", \"42\u00d740=\"],\n  [\"38\u00d762=\", \"38\u00d794=\"],\n  [\"83\u00d736=\", \"25\u00d721=\"],\n  [\"42\u00d778=\", \"36\u00d768=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the text replacements described by the diff:\n# - update the date line\n# - update each two-digit multiplication problem in the table\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2026-02-23 Monday\", \"2026-02-24 Tuesday\"),\n    @(\"99\u00d798=\", \"26\u00d738=\"),\n    @(\"43\u00d722=\", \"97\u00d793=\"),\n    @(\"67\u00d789=\", \"66\u00d796=\"),\n    @(\"51\u00d728=\", \"98\u00d747=\"),\n    @(\"18\u00d714=\", \"43\u00d726=\"),\n    @(\"47\u00d741=\", \"67\u00d768=\"),\n    @(\"60\u00d766=\", \"73\u00d734=\"),\n    @(\"69\u00d789=\", \"53\u00d740=\"),\n    @(\"17\u00d739=\", \"71\u00d721=\"),\n    @(\"91\u00d765=\", \"43\u00d797=\"),\n    @(\"41\u00d718=\", \"90\u00d773=\"),\n    @(\"35\u00d764=\", \"11\u00d734=\"),\n    @(\"19\u00d787=\", \"70\u00d782=\"),\n    @(\"58\u00d739=\", \"26\u00d714=\"),\n    @(\"21\u00d765=\", \"85\u00d798=\"),\n    @(\"43\u00d729=\", \"71\u00d711=\"),\n    @(\"70\u00d779=\", \"42\u00d777=\"),\n    @(\"73\u00d761=\", \"56\u00d770=\"),\n    @(\"40\u00d782=\", \"47\u00d743=\"),\n    @(\"92\u00d711=\", \"55\u00d743=\"),\n    @(\"96\u00d786=\", \"94\u00d773=\"),\n    @(\"70\u00d759=\", \"42\u00d740=\"),\n    @(\"38\u00d762=\", \"38\u00d794=\"),\n    @(\"83\u00d736=\", \"25\u00d721=\"),\n    @(\"42\u00d778=\", \"36\u00d768=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $range = $d.Content\n    $range.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
